# CS133JS Lab02 Rubrics - "Fixed typo and formatting"
#
# 1) Fix the typo "concole" -> "console" in the rubric line that appears
#    in all three sheets (Group A, Group B, Group C), and clear the
#    leftover stray font override that had been applied to that cell.
# 2) Update the selection in every sheet to A5.
# 3) Move the active/selected tab from "Group A" to "Group C".
# 4) Reset the zoom on "Group A" back to the (unzoomed) normal view.
# 5) Turn on portrait page orientation for "Group A" (page setup).

$wb = $excel.ActiveWorkbook

$fixedText = "Function exercises in the console"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A5")
    $cell.Value = $fixedText
    # The cell had picked up a stray "Calibri (Body)" font override;
    # drop back to the sheet's normal style.
    $cell.Style = "Normal"

    # Park the visible selection on A5, matching the saved view state.
    $ws.Range("A5").Select() | Out-Null
}

# "Group C" (3rd tab) becomes the active/selected sheet.
$groupC = $wb.Worksheets.Item("Group C")
$groupC.Activate() | Out-Null
$groupC.Range("A5").Select() | Out-Null

# "Group A" view goes back to its normal (unzoomed) view.
$groupA = $wb.Worksheets.Item("Group A")
$groupA.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 100
$groupA.Range("A5").Select() | Out-Null

# "Group A" page setup switches to portrait orientation.
$groupA.PageSetup.Orientation = 1

# Leave "Group C" as the final active sheet/tab.
$groupC.Activate() | Out-Null
$groupC.Range("A5").Select() | Out-Null
